$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.872.62'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '1.784.43'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Value = '''226.27'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = '''0.555'
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').Value = '''30.58'
$ws.Range('E8').Value = '  -4.40%  '
$ws.Range('D9').Value = '''46.32'
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('D10').Value = '''0.280'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '''0.0664'
$ws.Range('E11').Value = '  -1.48%  '
$ws.Range('D12').Value = '''0.0928'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '2.040.75'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '1.787.17'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = '''10.74'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').Value = '''0.626'
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').Value = '33.875.96'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('E18').Value = '  -3.59%  '
$ws.Range('D19').Value = '''69.17'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = '''252.36'
$ws.Range('E20').Value = '  -4.73%  '
$ws.Range('D21').Value = '0.0₃0741'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').Value = '''10.30'
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('E24').Value = '  -3.97%  '
$ws.Range('D25').Value = '''2.14'
$ws.Range('E25').Value = '  -1.79%  '
$ws.Range('D26').Value = '''158.56'
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').Value = '''16.52'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('D29').Value = '''6.96'
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').Value = '''3.83'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '''0.0516'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('D34').Value = '''3.62'
$ws.Range('E34').Value = '  +2.95%  '
$ws.Range('D35').Value = '''1.84'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').Value = '1.502.15'
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('D38').Value = '''0.626'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').Value = '''83.74'
$ws.Range('E39').Value = '  -3.37%  '
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '''2.83'
$ws.Range('E41').Value = '  +2.07%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '''2.35'
$ws.Range('E42').Value = '  +1.71%  '
$ws.Range('D43').Value = '''0.907'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '''0.0517'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '''2.03'
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '1.935.31'
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D48').Value = '''5.71'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').Value = '''11.66'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('D51').Value = '''50.91'
$ws.Range('E51').Value = '  -6.29%  '
